# [FIX] update data formatting for consistency
#
# The "income" sheet header row (row 1) labels several sub-total columns
# without the "Total" qualifier, while the grand-total column is labelled
# redundantly ("Total Liabilities and Equity"). Rename the headers so the
# terminology is applied consistently:
#   B1: "Current Assets"              -> "Total Current Assets"
#   E1: "Current Liabilities"         -> "Total Current Liabilities"
#   H1: "Equity"                      -> "Total Equity"
#   N1: "Total Liabilities and Equity"-> "Liabilities and Equity"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("income")

$ws.Range("B1").Value = "Total Current Assets"
$ws.Range("E1").Value = "Total Current Liabilities"
$ws.Range("H1").Value = "Total Equity"
$ws.Range("N1").Value = "Liabilities and Equity"

# Leave the cursor where the author's last save left it.
$ws.Range("O2").Select()
